$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 855, shifting existing rows 855-896 down to 856-897.
$ws.Rows.Item(855).Insert()

# Make sure the date-like text in column A is stored as plain text,
# matching how the rest of the date column is stored (not auto-converted
# into a date serial number by Excel's type inference).
$ws.Cells.Item(855,1).NumberFormat = "@"
$ws.Cells.Item(855,1).Value = "2026/02/23"
# Drop the temporary text-number-format so the cell's style matches its
# neighbours (no leftover explicit style index).
$ws.Cells.Item(855,1).ClearFormats()

$ws.Cells.Item(855,2).Value = "月"
$ws.Cells.Item(855,3).Value = 8
$ws.Cells.Item(855,4).Value = 201
